$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Remove B2, shift remaining row-2 data
$ws.Range("B2").Value = $null

# Update row 2 values
$ws.Range("C2").Value = 3.8883659508071853
$ws.Range("D2").Value = 0.16359730644805667
$ws.Range("E2").Value = 3.4270214280954638

# Update row 3 values
$ws.Range("B3").Value = 0.28617109816574898
$ws.Range("C3").Value = 6.5318022016907191
$ws.Range("D3").Value = 0.98312787159949644
$ws.Range("E3").Value = 8.0593215041399819

# Update selection
$ws.Range("B1:E3").Select()
